$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53715ad08294096f446025a82d1c6680dd843c51/e2e/a4b82b49-f1f5-468d-951b-ab75af631667.md"
$handbackDisplay = "a4b82b49-f1f5-468d-951b-ab75af631667.md"
$newStatus = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" everywhere ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- zh-cn: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$zhcn.Range("J2").Value = $zhcn.Range("G2").Value2
$zhcn.Range("K2").Value = "2016-09-04 09:03:06"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $handbackUrl, "", "", $handbackDisplay) | Out-Null
$zhcn.Range("I2").Style = $zhcn.Range("A2").Style

# --- de-de: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$dede.Range("J2").Value = $dede.Range("G2").Value2
$dede.Range("K2").Value = "2016-09-04 09:03:15"

$dede.Hyperlinks.Add($dede.Range("I2"), $handbackUrl, "", "", $handbackDisplay) | Out-Null
$dede.Range("I2").Style = $dede.Range("A2").Style

# --- Column width changes (wider columns for handback info) ---
$overview.Range("E:F").ColumnWidth = 29.9777047293527

$zhcn.Range("C:C").ColumnWidth = 29.9777047293527
$zhcn.Range("I:J").ColumnWidth = 40

$dede.Range("C:C").ColumnWidth = 29.9777047293527
$dede.Range("I:J").ColumnWidth = 40
